$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "CardPos", same format as E1 (bold header style)
$ws.Range("F1").Value = "CardPos"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New data cell F2 = 0, formatted as an integer number
$ws.Range("F2").Value = 0
$ws.Range("F2").NumberFormat = "0"
